# Test Suite Execution Engine
# Adds a "Page Objects" column + a new "Click LogOut Button" test step to
# the "Test Steps" sheet, and adds a new "Test Cases" sheet describing the
# test-case/run-mode matrix.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

# ---- Sheet 1 ("Test Steps"): final 8x5 grid, written cell-by-cell ----
# Row 1 (headers)
$ws1.Range("A1").Value = "Test Case ID"
$ws1.Range("B1").Value = "TS_ID"
$ws1.Range("C1").Value = "Description"
$ws1.Range("D1").Value = "Page Objects"
$ws1.Range("E1").Value = "Action_Keyword"

# Row 2
$ws1.Range("A2").Value = "Login_01"
$ws1.Range("B2").Value = "TS_001"
$ws1.Range("C2").Value = "Open Browser"
$ws1.Range("D2").ClearContents()
$ws1.Range("E2").Value = "openBrowser"

# Row 3
$ws1.Range("A3").Value = "Login_01"
$ws1.Range("B3").Value = "TS_002"
$ws1.Range("C3").Value = "Navigate to Website"
$ws1.Range("D3").ClearContents()
$ws1.Range("E3").Value = "navigateWebsite"

# Row 4
$ws1.Range("A4").Value = "Login_01"
$ws1.Range("B4").Value = "TS_003"
$ws1.Range("C4").Value = "Enter Username"
$ws1.Range("D4").Value = "txtbx_UserName"
$ws1.Range("E4").Value = "inputUsername"

# Row 5
$ws1.Range("A5").Value = "Login_01"
$ws1.Range("B5").Value = "TS_004"
$ws1.Range("C5").Value = "Enter Password"
$ws1.Range("D5").Value = "txtbx_Password"
$ws1.Range("E5").Value = "inputPassword"

# Row 6
$ws1.Range("A6").Value = "Login_01"
$ws1.Range("B6").Value = "TS_005"
$ws1.Range("C6").Value = "Click Login Button"
$ws1.Range("D6").Value = "btn_LogIn"
$ws1.Range("E6").Value = "doLogin"

# Row 7 (new step: Click LogOut Button)
$ws1.Range("A7").Value = "Login_01"
$ws1.Range("B7").Value = "TS_006"
$ws1.Range("C7").Value = "Click LogOut Button"
$ws1.Range("D7").Value = "btn_Signout"
$ws1.Range("E7").Value = "doLogout"

# Row 8 (old "Quit Browser" row, pushed down, renumbered TS_007)
$ws1.Range("A8").Value = "Login_01"
$ws1.Range("B8").Value = "TS_007"
$ws1.Range("C8").Value = "Quit Browser"
$ws1.Range("D8").ClearContents()
$ws1.Range("E8").Value = "doClose"

# ---- widen the new "Page Objects" column's neighbour (shifted Action_Keyword column) ----
$ws1.Columns.Item(5).ColumnWidth = 15.333333333333334

$ws1.Range("D7").Select()

# ---- Sheet 2 ("Test Cases") ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Test Cases"

$ws2.Range("A1").Value = "Test Case ID"
$ws2.Range("B1").Value = "Description"
$ws2.Range("C1").Value = "Runmode"

$ws2.Range("A2").Value = "Login_01"
$ws2.Range("B2").Value = "Login to the online app"
$ws2.Range("C2").Value = "Yes"

$ws2.Range("A3").Value = "Login_02"
$ws2.Range("B3").Value = "Login to the online app"
$ws2.Range("C3").Value = "No"

$ws2.Columns.Item(1).ColumnWidth = 11.5
$ws2.Columns.Item(2).ColumnWidth = 21.5

$ws2.Range("D8").Select()

# ---- restore "Test Steps" as the active/selected sheet+cell ----
$ws1.Activate()
$ws1.Range("D7").Select()
